$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab
$ws.Name = "Gamma1F"

# Append a new row of averaged intensity data (row 16), following the
# same pattern/formatting used by the previous data rows: column A
# carries the bold/bordered "index" style used throughout the table.
$ws.Cells.Item(15, 1).Copy($ws.Cells.Item(16, 1))

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "HexGrid-60degTilt5degRes"

$ws.Cells.Item(16, 3).Value = 1.013849313885671
$ws.Cells.Item(16, 4).Value = 0.9229541462827703
$ws.Cells.Item(16, 5).Value = 1.011873886281002
$ws.Cells.Item(16, 6).Value = 1.013849313885671
$ws.Cells.Item(16, 7).Value = 0.960652774054112
$ws.Cells.Item(16, 8).Value = 1.041031941877974
$ws.Cells.Item(16, 9).Value = 1.012512863396901
$ws.Cells.Item(16, 10).Value = 0.9229541462827703
$ws.Cells.Item(16, 11).Value = 0.9674140162818861
$ws.Cells.Item(16, 12).Value = 0.9906316650837783
$ws.Cells.Item(16, 13).Value = 0.9938124876297384
